$p = $ppt.ActivePresentation
$cxp = $p.CustomXMLParts
Get-Member -InputObject $cxp | Out-String | Write-Output
Write-Output "Count: $($cxp.Count)"
